$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N4").Value = 2.5
$ws.Range("O4").Value = 1.53
$ws.Range("J7").Value = 1.18
$ws.Range("K7").Value = 4.5
$ws.Range("R7").Value = 2.63
$ws.Range("S7").Value = 1.44
$ws.Range("W7").Value = 21
$ws.Range("AG7").Value = 17
$ws.Range("J9").Value = 1.06
$ws.Range("K9").Value = 10
$ws.Range("G11").Value = 2
$ws.Range("G16").Value = 3.6
$ws.Range("I16").Value = 1.9
$ws.Range("AF16").Value = 10
$ws.Range("G25").Value = 4.05
$ws.Range("H25").Value = 3.3
$ws.Range("I25").Value = 1.85
$ws.Range("L25").Value = 1.39
$ws.Range("M25").Value = 2.57
$ws.Range("N25").Value = 2.12
$ws.Range("R25").Value = 1.98
$ws.Range("S25").Value = 1.65
$ws.Range("T25").Value = 9.75
$ws.Range("U25").Value = 21
$ws.Range("V25").Value = 14
$ws.Range("W25").Value = 65
$ws.Range("X25").Value = 45
$ws.Range("AA25").Value = 6.5
$ws.Range("AB25").Value = 18
$ws.Range("AC25").Value = 110
$ws.Range("AE25").Value = 5.8
$ws.Range("AF25").Value = 7.7
$ws.Range("AH25").Value = 15
$ws.Range("H26").Value = 3.1
$ws.Range("I26").Value = 2.67
$ws.Range("P26").Value = 1.52
$ws.Range("Q26").Value = 2.22
$ws.Range("R26").Value = 1.98
$ws.Range("S26").Value = 1.65
$ws.Range("T26").Value = 6.6
$ws.Range("U26").Value = 11.25
$ws.Range("V26").Value = 10.25
$ws.Range("W26").Value = 27
$ws.Range("X26").Value = 25
$ws.Range("Y26").Value = 45
$ws.Range("Z26").Value = 7.2
$ws.Range("AA26").Value = 6.1
$ws.Range("AB26").Value = 18
$ws.Range("AC26").Value = 110
$ws.Range("AE26").Value = 6.8
$ws.Range("AF26").Value = 11.75
$ws.Range("AH26").Value = 30
$ws.Range("AI26").Value = 27
$ws.Range("G32").Value = 3.55
$ws.Range("I32").Value = 2.25
$ws.Range("J32").Value = 1.1
$ws.Range("K32").Value = 6
$ws.Range("L32").Value = 1.39
$ws.Range("M32").Value = 2.77
$ws.Range("N32").Value = 2.15
$ws.Range("O32").Value = 1.62
$ws.Range("Q32").Value = 2.65
$ws.Range("R32").Value = 1.8
$ws.Range("S32").Value = 1.91
$ws.Range("T32").Value = 9
$ws.Range("U32").Value = 19
$ws.Range("V32").Value = 11.75
$ws.Range("W32").Value = 55
$ws.Range("Z32").Value = 6
$ws.Range("AB32").Value = 13.5
$ws.Range("AC32").Value = 65
$ws.Range("AD32").Value = 600
$ws.Range("AE32").Value = 6.9
$ws.Range("AG32").Value = 8.75
$ws.Range("AH32").Value = 24
$ws.Range("AI32").Value = 19
$ws.Range("AJ32").Value = 29
$ws.Range("G33").Value = 3.55
$ws.Range("H33").Value = 2.95
$ws.Range("I33").Value = 2.18
$ws.Range("L33").Value = 1.45
$ws.Range("N33").Value = 2.32
$ws.Range("T33").Value = 8.5
$ws.Range("U33").Value = 18
$ws.Range("V33").Value = 12
$ws.Range("W33").Value = 50
$ws.Range("AE33").Value = 6
$ws.Range("AF33").Value = 9.25
$ws.Range("AG33").Value = 9.25
$ws.Range("AH33").Value = 21
$ws.Range("AI33").Value = 21
$ws.Range("G37").Value = 2.75
$ws.Range("I37").Value = 2.3
$ws.Range("K37").Value = 9
$ws.Range("L37").Value = 1.18
$ws.Range("M37").Value = 4.35
$ws.Range("N37").Value = 1.53
$ws.Range("O37").Value = 2.32
$ws.Range("P37").Value = 1.29
$ws.Range("Q37").Value = 3.3
$ws.Range("R37").Value = 1.47
$ws.Range("S37").Value = 2.52
$ws.Range("T37").Value = 12.5
$ws.Range("U37").Value = 17
$ws.Range("V37").Value = 10
$ws.Range("X37").Value = 20
$ws.Range("Y37").Value = 22
$ws.Range("Z37").Value = 9
$ws.Range("AA37").Value = 7.3
$ws.Range("AB37").Value = 11
$ws.Range("AC37").Value = 35
$ws.Range("AE37").Value = 12
$ws.Range("AF37").Value = 15
$ws.Range("AH37").Value = 26
$ws.Range("AI37").Value = 16
$ws.Range("AJ37").Value = 19
$ws.Range("G44").Value = 1.6
$ws.Range("H44").Value = 3.6
$ws.Range("I44").Value = 6.25
$ws.Range("J44").Value = 1.08
$ws.Range("K44").Value = 8
$ws.Range("R44").Value = 2.05
$ws.Range("S44").Value = 1.7
$ws.Range("T44").Value = 6
$ws.Range("U44").Value = 7
$ws.Range("AB44").Value = 19
$ws.Range("AG44").Value = 19
$ws.Range("L51").Value = 1.57
$ws.Range("M51").Value = 2.26
$ws.Range("O51").Value = 1.41
$ws.Range("P51").Value = 1.63
$ws.Range("Q51").Value = 2.17
$ws.Range("R51").Value = 2.3
$ws.Range("S51").Value = 1.55
$ws.Range("T51").Value = 4.2
$ws.Range("U51").Value = 7
$ws.Range("Z51").Value = 4.9
$ws.Range("AE51").Value = 5.6
$ws.Range("K64").Value = 9
$ws.Range("K67").Value = 15
$ws.Range("L67").Value = 1.2
$ws.Range("M67").Value = 4.33
$ws.Range("N67").Value = 1.65
$ws.Range("O67").Value = 2.2
$ws.Range("R67").Value = 1.91
$ws.Range("S67").Value = 1.91
$ws.Range("U67").Value = 7
$ws.Range("J70").Value = 1.05
$ws.Range("K70").Value = 11
$ws.Range("L70").Value = 1.29
$ws.Range("M70").Value = 3.5
$ws.Range("N70").Value = 1.93
$ws.Range("O70").Value = 1.93
$ws.Range("K71").Value = 17
$ws.Range("N71").Value = 1.53
$ws.Range("O71").Value = 2.4
$ws.Range("H72").Value = 3.65
$ws.Range("Q72").Value = 2.95
$ws.Range("T72").Value = 8.25
$ws.Range("U72").Value = 9.75
$ws.Range("W72").Value = 16.5
$ws.Range("X72").Value = 14.5
$ws.Range("AA72").Value = 7.1
$ws.Range("AE72").Value = 12
$ws.Range("AI72").Value = 29
$ws.Range("AJ72").Value = 32
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 3.65
$ws.Range("I73").Value = 3.3
$ws.Range("K73").Value = 8.25
$ws.Range("L73").Value = 1.23
$ws.Range("M73").Value = 3.75
$ws.Range("N73").Value = 1.7
$ws.Range("O73").Value = 2.05
$ws.Range("P73").Value = 1.34
$ws.Range("Q73").Value = 3
$ws.Range("R73").Value = 1.62
$ws.Range("S73").Value = 2.18
$ws.Range("T73").Value = 9
$ws.Range("U73").Value = 10.75
$ws.Range("V73").Value = 8.5
$ws.Range("W73").Value = 18.5
$ws.Range("X73").Value = 14.5
$ws.Range("Z73").Value = 8.25
$ws.Range("AA73").Value = 7.1
$ws.Range("AB73").Value = 13
$ws.Range("AC73").Value = 50
$ws.Range("AD73").Value = 300
$ws.Range("AE73").Value = 11.75
$ws.Range("AF73").Value = 18.5
$ws.Range("AG73").Value = 11.5
$ws.Range("AH73").Value = 45
$ws.Range("AI73").Value = 27
$ws.Range("AJ73").Value = 30
